$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.040.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.049.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.79"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.04"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +8.38%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0816"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.71"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.353.50"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.02"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.758"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.27"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.045.15"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.989.26"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.68"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0826"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.53"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.19"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.18"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.94"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.47"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0603"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.29"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +9.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.30"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.23"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.533.39"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.55"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.84"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.60"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0925"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.21%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.05"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.242.63"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.88%  "
